$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.776.41'
$ws.Range('E2').Value = '  -3.49%  '
$ws.Range('D3').Value = '3.478.35'
$ws.Range('E3').Value = '  -3.11%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '604.86'
$ws.Range('E5').Value = '  -3.32%  '
$ws.Range('D6').Value = '148.83'
$ws.Range('E6').Value = '  -5.78%  '
$ws.Range('D7').Value = '3.476.64'
$ws.Range('E7').Value = '  -3.11%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -2.16%  '
$ws.Range('D10').Value = '0.143'
$ws.Range('E10').Value = '  -3.98%  '
$ws.Range('E11').Value = '  +1.91%  '
$ws.Range('E12').Value = '  -4.16%  '
$ws.Range('E13').Value = '  -5.16%  '
$ws.Range('D14').Value = '31.76'
$ws.Range('E14').Value = '  -5.76%  '
$ws.Range('D15').Value = '4.065.48'
$ws.Range('E15').Value = '  -3.09%  '
$ws.Range('D16').Value = '3.477.31'
$ws.Range('E16').Value = '  -2.74%  '
$ws.Range('D17').Value = '66.837.33'
$ws.Range('E17').Value = '  -3.79%  '
$ws.Range('E18').Value = '  -0.60%  '
$ws.Range('D19').Value = '6.48'
$ws.Range('E19').Value = '  -6.02%  '
$ws.Range('D20').Value = '15.42'
$ws.Range('E20').Value = '  -4.67%  '
$ws.Range('D21').Value = '10.14'
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('D22').Value = '440.45'
$ws.Range('E22').Value = '  -4.79%  '
$ws.Range('D23').Value = '0.614'
$ws.Range('E23').Value = '  -5.49%  '
$ws.Range('D24').Value = '79.44'
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('D26').Value = '3.614.09'
$ws.Range('E26').Value = '  -2.96%  '
$ws.Range('D27').Value = '0.0000121'
$ws.Range('E27').Value = '  -9.73%  '
$ws.Range('D28').Value = '9.80'
$ws.Range('E28').Value = '  -8.47%  '
$ws.Range('D29').Value = '8.43'
$ws.Range('E29').Value = '  -8.97%  '
$ws.Range('E30').Value = '  -4.24%  '
$ws.Range('E31').Value = '  -7.26%  '
$ws.Range('D32').Value = '0.169'
$ws.Range('E32').Value = '  -2.05%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').Value = '25.47'
$ws.Range('E34').Value = '  -3.70%  '
$ws.Range('D35').Value = '6.08'
$ws.Range('E35').Value = '  -7.43%  '
$ws.Range('D36').Value = '3.465.98'
$ws.Range('E36').Value = '  -3.22%  '
$ws.Range('E37').Value = '  -7.39%  '
$ws.Range('D38').Value = '7.95'
$ws.Range('E38').Value = '  -5.50%  '
$ws.Range('E40').Value = '  +0.25%  '
$ws.Range('D41').Value = '176.18'
$ws.Range('E41').Value = '  -1.88%  '
$ws.Range('D42').Value = '0.0892'
$ws.Range('E42').Value = '  -3.85%  '
$ws.Range('D43').Value = '2.15'
$ws.Range('E43').Value = '  -11.33%  '
$ws.Range('E44').Value = '  -4.45%  '
$ws.Range('D45').Value = '0.890'
$ws.Range('E45').Value = '  -1.97%  '
$ws.Range('D46').Value = '29.26'
$ws.Range('E46').Value = '  -6.53%  '
$ws.Range('D47').Value = '46.31'
$ws.Range('E47').Value = '  +0.63%  '
$ws.Range('D48').Value = '1.25'
$ws.Range('E48').Value = '  -9.52%  '
$ws.Range('E49').Value = '  -4.93%  '
$ws.Range('E50').Value = '  -9.60%  '
$ws.Range('D51').Value = '0.987'
$ws.Range('E51').Value = '  -5.04%  '
